$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "is_active" column (E) currently holds =TRUE()/=FALSE() formulas that
# evaluate to numeric booleans. Replace them with literal text values
# "TRUE"/"FALSE" (stored as plain text, not booleans), matching rows:
#   2,3,4,5,6,7,10,11 -> "TRUE"
#   8,9               -> "FALSE"
#
# A direct $cell.Value = "TRUE" assignment gets auto-coerced by Excel into a
# Boolean cell (t="b"), so instead we build the text via a formula that
# evaluates to a text string, then convert that formula result to a
# literal value in-place (paste values only), which keeps the cell's
# existing style/number format untouched and stores a genuine text value.

$trueRows = @(2,3,4,5,6,7,10,11)
$falseRows = @(8,9)

foreach ($r in $trueRows) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Formula = "=""TRUE"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

foreach ($r in $falseRows) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Formula = "=""FALSE"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0

# Move the active selection to F11, as recorded in the saved view state.
[void]$ws.Range("F11").Select()
